# Applies the betexplorer costa-rica primera-division 2023-2024 update:
#  - rows 38/39 swap their match data (F:V)
#  - rows 91/92 swap their match data (F:V)
#  - six new match rows (97-102) are appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Rows 38 <-> 39 and 91 <-> 92: swap F:V (A:E stay tied to the
#    physical row - they hold index/country/tournament/season/date).
# ---------------------------------------------------------------
function Swap-Row($r1, $r2) {
    for ($col = 6; $col -le 22; $col++) {
        $v1 = $ws.Cells.Item($r1, $col).Value2
        $v2 = $ws.Cells.Item($r2, $col).Value2
        $ws.Cells.Item($r1, $col).Value = $v2
        $ws.Cells.Item($r2, $col).Value = $v1
    }
}

Swap-Row 38 39
Swap-Row 91 92

# ---------------------------------------------------------------
# 2) Append six new rows (97-102), copying the style of row 96's
#    A (index) and E (date) cells, then filling in the values.
# ---------------------------------------------------------------
$ws.Range("A96").Copy()
$ws.Range("A97:A102").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E96").Copy()
$ws.Range("E97:E102").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$newRows = @(
    @{ A=96;  E=45227.125;          F="AD Santos";     G=1; H="Cartagines";        I=0; J=2.73; K="25/10/2023 13:42"; L=2.54; M="28/10/2023 02:51"; N=3.55; O="25/10/2023 13:42"; P=3.33; Q="28/10/2023 02:50"; R=2.37; S="25/10/2023 13:42"; T=2.86; U="28/10/2023 02:51"; V="https://www.betexplorer.com/football/costa-rica/primera-division/santos-de-guapiles-cartagines/8EeIsBRO/" },
    @{ A=97;  E=45227.16666666666;  F="Liberia";       G=2; H="Sporting San Jose"; I=1; J=2.03; K="24/10/2023 05:42"; L=2.21; M="28/10/2023 03:59"; N=3.59; O="24/10/2023 05:42"; P=3.76; Q="28/10/2023 03:50"; R=3.52; S="24/10/2023 05:42"; T=3.06; U="28/10/2023 03:59"; V="https://www.betexplorer.com/football/costa-rica/primera-division/liberia-sporting-san-jose/t2QW0YRn/" },
    @{ A=98;  E=45227.95833333334;  F="Puntarenas FC"; G=1; H="Guanacasteca";      I=2; J=2.13; K="25/10/2023 11:43"; L=2.23; M="28/10/2023 22:51"; N=3.28; O="25/10/2023 11:43"; P=3.29; Q="28/10/2023 22:57"; R=3.37; S="25/10/2023 11:43"; T=3.41; U="28/10/2023 22:57"; V="https://www.betexplorer.com/football/costa-rica/primera-division/puntarenas-fc-guanacasteca/04tYbCda/" },
    @{ A=99;  E=45228.04166666666;  F="Zeledon";       G=0; H="Grecia";            I=2; J=1.75; K="25/10/2023 11:43"; L=1.85; M="29/10/2023 00:52"; N=3.9;  O="25/10/2023 11:43"; P=3.77; Q="29/10/2023 00:57"; R=4.03; S="25/10/2023 11:43"; T=4.14; U="29/10/2023 00:52"; V="https://www.betexplorer.com/football/costa-rica/primera-division/zeledon-grecia/KhPz0hsg/" },
    @{ A=100; E=45228.125;          F="Alajuelense";   G=1; H="Herediano";         I=0; J=1.93; K="25/10/2023 11:44"; L=2.04; M="29/10/2023 02:51"; N=3.64; O="25/10/2023 11:44"; P=3.77; Q="29/10/2023 02:51"; R=3.55; S="25/10/2023 11:44"; T=3.45; U="29/10/2023 02:51"; V="https://www.betexplorer.com/football/costa-rica/primera-division/alajuelense-herediano/nuRS1ECt/" },
    @{ A=101; E=45228.95833333334;  F="Saprissa";      G=2; H="San Carlos";        I=1; J=1.28; K="25/10/2023 11:44"; L=1.42; M="29/10/2023 22:57"; N=5.38; O="25/10/2023 11:44"; P=4.96; Q="29/10/2023 22:57"; R=8.58; S="25/10/2023 11:44"; T=6.87; U="29/10/2023 22:57"; V="https://www.betexplorer.com/football/costa-rica/primera-division/saprissa-san-carlos/AJZ06f4P/" }
)

$r = 97
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value  = $row.A
    $ws.Cells.Item($r, 2).Value  = "costa-rica"
    $ws.Cells.Item($r, 3).Value  = "primera-division"
    $ws.Cells.Item($r, 4).Value  = "2023-2024"
    $ws.Cells.Item($r, 5).Value  = $row.E
    $ws.Cells.Item($r, 6).Value  = $row.F
    $ws.Cells.Item($r, 7).Value  = $row.G
    $ws.Cells.Item($r, 8).Value  = $row.H
    $ws.Cells.Item($r, 9).Value  = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $ws.Cells.Item($r, 21).Value = $row.U
    $ws.Cells.Item($r, 22).Value = $row.V
    $r++
}
